# 10th - MarketBeat update for CURO (single stock) - adds a new
# "Jun_27" week column (replacing the unused "Jun_15" column) and
# appends two new analyst firms to the bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Drop the empty "Jun_15" column (column C). This shifts the
#    "Jun_13" column (old D) into C, and the "Jun_10" column (old E)
#    into D - matching the data that is still present after the edit.
# ------------------------------------------------------------------
$ws.Columns("C").Delete()

# ------------------------------------------------------------------
# 2) The old "Jun_17" column (B) becomes the new "Jun_27" week,
#    carrying the new downgrade entry for Zacks Investment Research
#    (row 5), highlighted the same way the other "new" entries are.
# ------------------------------------------------------------------
$ws.Range("B1").Value = "Jun_27"
$ws.Range("B5").Value = "6/27/2018,Downgrades,Buy -> Hold,"
$ws.Range("B5").Interior.ColorIndex = 45
$ws.Range("B5").Interior.Pattern = -4142

# ------------------------------------------------------------------
# 3) Append two newly tracked research firms at the bottom of the
#    table (no ratings yet, so "UN" like every other firm started
#    with).
# ------------------------------------------------------------------
$ws.Cells.Item(28, 1).Value = "Benchmark"
$ws.Cells.Item(28, 2).Value = "UN"
$ws.Cells.Item(29, 1).Value = "Evercore ISI"
$ws.Cells.Item(29, 2).Value = "UN"
